$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values changed for columns B-E
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 data values changed for columns B-E
$ws.Range("B2").Value = 10.711959559265367
$ws.Range("C2").Value = 15.412007316643699
$ws.Range("D2").Value = 8.2976607373479681
$ws.Range("E2").Value = 13.890536628298698

# Row 3 data values changed for columns B-E
$ws.Range("B3").Value = 13.314304289333535
$ws.Range("C3").Value = 15.463786012346466
$ws.Range("D3").Value = 14.297854565260399
$ws.Range("E3").Value = 15.294903384812345

# Selection changed from B1:AY3 to B1:E3
$ws.Range("B1:E3").Select()
